$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "482"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2631981.26"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "65"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "873168.55"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "138"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "790538.38"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "530"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2687400.85"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "243"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2563279.45"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "409"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2001418.71"

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "168"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1410508.28"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "9"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "272000.00"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "379"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3222413.95"

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "144"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1852274.71"

$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "15"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "363767.27"

$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "3113"
$ws.Range("D64").NumberFormat = "@"
$ws.Range("D64").Value = "18558487.38"

$ws.Range("C90").NumberFormat = "@"
$ws.Range("C90").Value = "288"
$ws.Range("D90").NumberFormat = "@"
$ws.Range("D90").Value = "1303356.14"

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "1098"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "5768914.79"

$ws.Range("C92").NumberFormat = "@"
$ws.Range("C92").Value = "471"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "4317852.63"

$ws.Range("C105").NumberFormat = "@"
$ws.Range("C105").Value = "495"
$ws.Range("D105").NumberFormat = "@"
$ws.Range("D105").Value = "2306017.74"

$ws.Range("C106").NumberFormat = "@"
$ws.Range("C106").Value = "231"
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "1756097.68"

$ws.Range("C108").NumberFormat = "@"
$ws.Range("C108").Value = "18"
$ws.Range("D108").NumberFormat = "@"
$ws.Range("D108").Value = "357957.62"

$ws.Range("C111").NumberFormat = "@"
$ws.Range("C111").Value = "1769"
$ws.Range("D111").NumberFormat = "@"
$ws.Range("D111").Value = "7575412.91"

$ws.Range("C112").NumberFormat = "@"
$ws.Range("C112").Value = "737"
$ws.Range("D112").NumberFormat = "@"
$ws.Range("D112").Value = "5120130.94"

$ws.Range("C113").NumberFormat = "@"
$ws.Range("C113").Value = "254"
$ws.Range("D113").NumberFormat = "@"
$ws.Range("D113").Value = "3217212.51"

$ws.Range("C114").NumberFormat = "@"
$ws.Range("C114").Value = "84"
$ws.Range("D114").NumberFormat = "@"
$ws.Range("D114").Value = "1250060.00"

$ws.Range("C115").NumberFormat = "@"
$ws.Range("C115").Value = "12"
$ws.Range("D115").NumberFormat = "@"
$ws.Range("D115").Value = "356903.00"

